$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 loses its value but the "Yes/No" click leaves the row itself behind
# (touching a row-level property keeps an explicit, empty <row r="5"/> in
# the saved XML instead of the row disappearing entirely).
$ws.Rows(5).OutlineLevel = 0
$ws.Range("A5").ClearContents()

# These rows' only content was the A-column value; clearing it removes the
# row from the sheet data entirely.
$ws.Range("A8").ClearContents()
$ws.Range("A10").ClearContents()
$ws.Range("A20").ClearContents()
$ws.Range("A28").ClearContents()
$ws.Range("A29").ClearContents()

# Row 38 keeps its B38 value but loses A38.
$ws.Range("A38").ClearContents()
